# Explications de la fonction Player Bluetooth.docx
#
# 1) In the "Fonction 1" paragraph, fix the typo "fichiers sons" -> "fichiers son"
#    and split the surrounding sentence into three runs (mirroring the run
#    fragmentation Word performs around the edit point), with the
#    "_GoBack" bookmark now sitting right after "... fichiers son".
# 2) Remove the "_GoBack" bookmark that used to sit in the "fournis par
#    Carl" paragraph further down (Word only ever keeps one "_GoBack"
#    bookmark - it moves to the most recent edit location).

$d = $word.ActiveDocument

# --- Move the "_GoBack" bookmark -----------------------------------------

# Delete the old "_GoBack" bookmark (currently between "fournis par Carl "
# and "(liste non exhaustive)...") so we can re-add it elsewhere with the
# same reserved name.
$d.Bookmarks("_GoBack").Delete()

# --- Fix "Fonction 1" paragraph text + run layout -------------------------

# Locate the paragraph by its known text (more robust than a hard-coded
# paragraph index) and use its start as the offset anchor.
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$null = $anchor.Find.Execute("Pouvoir se connecter en Bluetooth", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pStart = $anchor.Start

# Within this paragraph (offsets relative to its start):
#   0                                                                    80
#   "Pouvoir se connecter en Bluetooth grâce à un smartphone sur le système. Grâce à "
#   80                                                                  111
#   "ce lien, jouer des fichiers son" (was "sons")
#   111 (new "_GoBack" bookmark goes here)
#   111...
#   " présents sur le smartphone à travers les enceintes du système. On
#    pourra assimiler notre système à une enceinte Bluetooth nomade."

# Step 1: delete the extra "s" that turns "sons" into "son".
$sCut = $pStart + 111
$sRange = $d.Range($sCut, $sCut + 1)
$sRange.Delete()

# Step 2: force a run split right before "ce lien" (offset 80) by adding
# and immediately removing a throwaway bookmark at that collapsed point -
# bookmarks can't live inside a run, so Word splits the run around them,
# and the split survives the bookmark's removal.
$splitCut = $pStart + 80
$splitRange = $d.Range($splitCut, $splitCut)
$d.Bookmarks.Add("ZZTempSplit", $splitRange) | Out-Null
$d.Bookmarks("ZZTempSplit").Delete()

# Step 3: re-add the "_GoBack" bookmark right after "... fichiers son",
# i.e. where the user's last edit (the dropped "s") happened.
$goBackCut = $pStart + 111
$goBackRange = $d.Range($goBackCut, $goBackCut)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
